$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, shifting existing rows 141:274 down to 142:275
$ws.Rows("141:141").Insert()

# Populate the new row 141 with fresh data
$ws.Range("A141").Value = 10
$ws.Range("B141").Value = "Vega Modelo de Temuco"
$ws.Range("C141").Value = "La Araucanía"
$ws.Range("D141").Value = 44484
$ws.Range("E141").Value = 9
$ws.Range("F141").Value = 100112023
$ws.Range("G141").Value = "Brócoli"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 500
$ws.Range("K141").Value = 700
$ws.Range("L141").Value = 800
$ws.Range("M141").Value = 760
$ws.Range("N141").Value = "$/unidad"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 760
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = "Hortaliza"
